$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (duplicate monte-carlo run), keep only row 1 (header) + row 2 (data)
$ws.Rows.Item(3).Delete()

# Re-state column R's (18th column) width explicitly -- matches the
# explicit <col .../> written for that column in the target workbook
$ws.Columns.Item(18).ColumnWidth = 8.3

# New header values
$ws.Range("T1").Value = "Standard Error"
$ws.Range("U1").Value = "Arrivals"
$ws.Range("V1").Value = "Iterations"
$ws.Range("W1").Value = "Lead Time Requirement"
$ws.Range("X1").Value = "Idle Time Requirement"

# Updated row 2 values
$ws.Range("Q2").Value = 0.16961226948758149
$ws.Range("R2").Value = 0.017034686872567876
$ws.Range("S2").Value = 0.33569177934307826
$ws.Range("T2").Value = 0.0045314641980172569
$ws.Range("U2").Value = 10000
$ws.Range("V2").Value = 5
$ws.Range("W2").Value = 0.5
$ws.Range("X2").Value = 50

# Move/select new used range
$ws.Range("A2:X2").Select() | Out-Null
